# Routing Table.xlsx edit
# Commit: "Add Item and Item's feature / Creating item@index function"
#
# Summary of changes applied to Sheet1:
#  - Rename several "nama method" values from singular user./profile. style
#    to plural users./profiles. style (H3, H4, H7, H9, H10)
#  - Rename a route path B10 from "Users/{Username}/Update" to
#    "Users/{Username}/edit"
#  - Rename route path B15 from "Items/{ItemID}" to "Items/{Username}/{ItemID}"
#  - Mark E17 (Items/create) as Private ("v")
#  - Fill in the previously empty "nama method" (G) / "nama view" (H) columns
#    for the whole Item section (rows 15-24), the Service section (rows 26-27),
#    the Category section (rows 29-31) and the Search row (row 35), applying
#    the same conditional-format "status" styles (Good/Neutral/Bad) already
#    used elsewhere on the sheet.
#  - Move the active selection/viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Plain value edits (style is unchanged on these cells)
# ---------------------------------------------------------------------------
$ws.Range("H3").Value  = "users.index"
$ws.Range("H4").Value  = "users.login"
$ws.Range("H7").Value  = "users.register"
$ws.Range("H9").Value  = "profiles.show"
$ws.Range("B10").Value = "Users/{Username}/edit"
$ws.Range("H10").Value = "profiles.edit"
$ws.Range("B15").Value = "Items/{Username}/{ItemID}"
$ws.Range("E17").Value = "v"

# ---------------------------------------------------------------------------
# 2. Helper: stamp a cell with one of the sheet's existing conditional
#    formatting "status" styles by copying format from a template cell that
#    already carries it, then set the new text. This reuses the existing
#    style indices (Good/Neutral/Bad) instead of minting new duplicate ones.
#
# NOTE: this COM-interop PowerShell host only binds POSITIONAL parameters
# reliably, so call this as: Set-StatusCell <template> <target> <value>
# ---------------------------------------------------------------------------
function Set-StatusCell {
    param(
        [string]$TemplateAddress,
        [string]$TargetAddress,
        [string]$NewValue
    )
    $ws.Range($TemplateAddress).Copy() | Out-Null
    $ws.Range($TargetAddress).PasteSpecial(-4122) | Out-Null
    if ($null -ne $NewValue) {
        $ws.Range($TargetAddress).Value = $NewValue
    }
}

# Template cells already on the sheet for each status style:
#   "Good"    (green)  center+center -> G3
#   "Neutral" (yellow) center+center -> G12
#   "Neutral" (yellow) no-alignment  -> H5
#   "Bad"     (red)    center+center -> H3
$Good         = "G3"
$Neutral      = "G12"
$NeutralBlank = "H5"
$Bad          = "H3"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Item section (rows 15-24)
# ---------------------------------------------------------------------------
Set-StatusCell $Bad     "G15" "item.show"
Set-StatusCell $Neutral "H15" "items.show"

Set-StatusCell $Good    "G16" "item.index"
Set-StatusCell $Neutral "H16" "items.index"

Set-StatusCell $Bad     "G17" "item.create"
Set-StatusCell $Neutral "H17" "items.create"

Set-StatusCell $Bad     "G18" "item.store"
Set-StatusCell $NeutralBlank "H18" $null

Set-StatusCell $Bad     "G19" "item.edit"
Set-StatusCell $Neutral "H19" "items.edit"

Set-StatusCell $Bad     "G20" "item.update"
Set-StatusCell $NeutralBlank "H20" $null

Set-StatusCell $Bad     "G21" "item.delete"
Set-StatusCell $Neutral "H21" "items.delete"

Set-StatusCell $Bad     "G22" "item.destroy"
Set-StatusCell $NeutralBlank "H22" $null

Set-StatusCell $Neutral "G23" "item.buy"
Set-StatusCell $Neutral "H23" "items.buy"

Set-StatusCell $Neutral "G24" "item.transaction.add"
Set-StatusCell $NeutralBlank "H24" $null

# ---------------------------------------------------------------------------
# 4. Services section (rows 26-27)
# ---------------------------------------------------------------------------
Set-StatusCell $Neutral "G26" "service.index"
Set-StatusCell $Neutral "H26" "services.index"

Set-StatusCell $Neutral "G27" "service.show"
Set-StatusCell $Neutral "H27" "services.show"

# ---------------------------------------------------------------------------
# 5. Category section (rows 29-31)
# ---------------------------------------------------------------------------
Set-StatusCell $Neutral "G29" "category.index"
Set-StatusCell $Neutral "G30" "category.show"
Set-StatusCell $Neutral "G31" "category.subcategory.show"

# ---------------------------------------------------------------------------
# 6. Search row (35) - style stays the plain default, only the value is new
# ---------------------------------------------------------------------------
$ws.Range("G35").Value = "search.engine"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Update the active selection / viewport to match the saved workbook state
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
